$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8491011808501696
$ws.Range("C2").Value = 0.04467913352849706
$ws.Range("D2").Value = 0.500737560853068
$ws.Range("E2").Value = 0.1574281767461958
$ws.Range("G2").Value = 0.002544653869563485
$ws.Range("I2").Value = 1.519282554190269
$ws.Range("J2").Value = 0.05803402912153022
$ws.Range("K2").Value = 0.5378763956628632
$ws.Range("L2").Value = 0.4346228824781164
$ws.Range("M2").Value = 0.2933735593118598
$ws.Range("N2").Value = 2.792936057651859
$ws.Range("O2").Value = 6.522228300717671
$ws.Range("B3").Value = 0.8186103464885832
$ws.Range("C3").Value = 0.04205291901339336
$ws.Range("D3").Value = 0.4993420319151767
$ws.Range("E3").Value = 0.1579034730284459
$ws.Range("G3").Value = 0.002547282046846016
$ws.Range("I3").Value = 1.527692481185589
$ws.Range("J3").Value = 0.05797511702315372
$ws.Range("K3").Value = 0.5075410874301269
$ws.Range("L3").Value = 0.4330334067178541
$ws.Range("M3").Value = 0.2878235080860456
$ws.Range("N3").Value = 2.815478194351568
$ws.Range("O3").Value = 6.551025228641151
$ws.Range("B4").Value = 0.8002274915249927
$ws.Range("C4").Value = 0.04042203794264765
$ws.Range("D4").Value = 0.4986789552701936
$ws.Range("E4").Value = 0.1582388360435232
$ws.Range("G4").Value = 0.00254898365746222
$ws.Range("I4").Value = 1.533392317450126
$ws.Range("J4").Value = 0.05793877849605966
$ws.Range("K4").Value = 0.4891135300305649
$ws.Range("L4").Value = 0.4322161356690515
$ws.Range("M4").Value = 0.284534262075038
$ws.Range("N4").Value = 2.830030724159098
$ws.Range("O4").Value = 6.571191189803159
$ws.Range("B5").Value = 0.7928221302171323
$ws.Range("C5").Value = 0.03975282724029938
$ws.Range("D5").Value = 0.4984576050388938
$ws.Range("E5").Value = 0.1583864718586465
$ws.Range("G5").Value = 0.002549699245461667
$ws.Range("I5").Value = 1.535850007852584
$ws.Range("J5").Value = 0.05792392875858798
$ws.Range("K5").Value = 0.4816545583144602
$ws.Range("L5").Value = 0.4319231105920807
$ws.Range("M5").Value = 0.2832238185839131
$ws.Range("N5").Value = 2.836140026388657
$ws.Range("O5").Value = 6.580034447044056
$ws.Range("B6").Value = 0.7915976768595385
$ws.Range("C6").Value = 0.03964142688354855
$ws.Range("D6").Value = 0.4984238047933758
$ws.Range("E6").Value = 0.158411650186542
$ws.Range("G6").Value = 0.00254981940928506
$ws.Range("I6").Value = 1.536266261264782
$ws.Range("J6").Value = 0.05792146047484437
$ws.Range("K6").Value = 0.4804190611623227
$ws.Range("L6").Value = 0.4318768746057771
$ws.Range("M6").Value = 0.2830080338546921
$ws.Range("N6").Value = 2.837165287177649
$ws.Range("O6").Value = 6.581540660573921
$ws.Range("B7").Value = 0.8001272719677957
$ws.Range("C7").Value = 0.04041303140369479
$ws.Range("D7").Value = 0.4986757720621711
$ws.Range("E7").Value = 0.1582407826479688
$ws.Range("G7").Value = 0.002548993218106967
$ws.Range("I7").Value = 1.533424916069492
$ws.Range("J7").Value = 0.05793857839507055
$ws.Range("K7").Value = 0.4890127309594021
$ws.Range("L7").Value = 0.4322120216333616
$ws.Range("M7").Value = 0.2845164675048295
$ws.Range("N7").Value = 2.830112391336757
$ws.Range("O7").Value = 6.571307919644624
$ws.Range("B8").Value = 0.8385180335322389
$ws.Range("C8").Value = 0.04377743199422213
$ws.Range("D8").Value = 0.5002162393951579
$ws.Range("E8").Value = 0.1575830404477436
$ws.Range("G8").Value = 0.002545541860525148
$ws.Range("I8").Value = 1.522071118853948
$ws.Range("J8").Value = 0.05801375090902283
$ws.Range("K8").Value = 0.5273758881635899
$ws.Range("L8").Value = 0.4340419690267154
$ws.Range("M8").Value = 0.2914353936967622
$ws.Range("N8").Value = 2.800560962215732
$ws.Range("O8").Value = 6.531642236888018
$ws.Range("B9").Value = 0.9164651693581334
$ws.Range("C9").Value = 0.05022918302607593
$ws.Range("D9").Value = 0.5047695650459332
$ws.Range("E9").Value = 0.1566374563727546
$ws.Range("G9").Value = 0.002539468234520856
$ws.Range("I9").Value = 1.504053701418307
$ws.Range("J9").Value = 0.05815983788212709
$ws.Range("K9").Value = 0.6041623270315597
$ws.Range("L9").Value = 0.4388847917768857
$ws.Range("M9").Value = 0.3059379930968902
$ws.Range("N9").Value = 2.748251384547864
$ws.Range("O9").Value = 6.473542938249977
$ws.Range("B10").Value = 0.9753306610190293
$ws.Range("C10").Value = 0.05488074863222892
$ws.Range("D10").Value = 0.5090431638544288
$ws.Range("E10").Value = 0.1561511317268156
$ws.Range("G10").Value = 0.002535425161159306
$ws.Range("I10").Value = 1.493397360837626
$ws.Range("J10").Value = 0.05826636055000467
$ws.Range("K10").Value = 0.6615078017956932
$ws.Range("L10").Value = 0.4432019760811556
$ws.Range("M10").Value = 0.317156539942431
$ws.Range("N10").Value = 2.713249779886038
$ws.Range("O10").Value = 6.44282303829965
$ws.Range("B11").Value = 1.002451977611742
$ws.Range("C11").Value = 0.05697774192096006
$ws.Range("D11").Value = 0.5111876934736728
$ws.Range("E11").Value = 0.1559748363220503
$ws.Range("G11").Value = 0.002533676009579202
$ws.Range("I11").Value = 1.489108272895898
$ws.Range("J11").Value = 0.05831464599870007
$ws.Range("K11").Value = 0.6877942559505925
$ws.Range("L11").Value = 0.4453297098727091
$ws.Range("M11").Value = 0.3223812581444321
$ws.Range("N11").Value = 2.698069681239517
$ws.Range("O11").Value = 6.431438947153424
$ws.Range("B12").Value = 1.012770883374913
$ws.Range("C12").Value = 0.05776908072789411
$ws.Range("D12").Value = 0.5120284768903076
$ws.Range("E12").Value = 0.1559145137060511
$ws.Range("G12").Value = 0.002533026535763577
$ws.Range("I12").Value = 1.487564285135562
$ws.Range("J12").Value = 0.05833290550359838
$ws.Range("K12").Value = 0.6977765352760343
$ws.Range("L12").Value = 0.4461588763886795
$ws.Range("M12").Value = 0.3243770379228366
$ws.Range("N12").Value = 2.692428020066391
$ws.Range("O12").Value = 6.427499986162758
$ws.Range("B13").Value = 1.010546369527219
$ws.Range("C13").Value = 0.05759877410855552
$ws.Range("D13").Value = 0.5118461248806909
$ws.Range("E13").Value = 0.1559272193892784
$ws.Range("G13").Value = 0.002533165839242521
$ws.Range("I13").Value = 1.487893245652998
$ws.Range("J13").Value = 0.05832897411665527
$ws.Range("K13").Value = 0.6956254287030674
$ws.Range("L13").Value = 0.4459792599241155
$ws.Range("M13").Value = 0.3239464444526732
$ws.Range("N13").Value = 2.69363830711022
$ws.Range("O13").Value = 6.428331778479901
$ws.Range("B14").Value = 1.003299949080031
$ws.Range("C14").Value = 0.0570429009584501
$ws.Range("D14").Value = 0.5112562908396541
$ws.Range("E14").Value = 0.1559697446881749
$ws.Range("G14").Value = 0.002533622319129952
$ws.Range("I14").Value = 1.488979641514888
$ws.Range("J14").Value = 0.05831614872509761
$ws.Range("K14").Value = 0.6886149426601662
$ws.Range("L14").Value = 0.4453974568051677
$ws.Range("M14").Value = 0.3225451064812646
$ws.Range("N14").Value = 2.69760339985864
$ws.Range("O14").Value = 6.431107435340039
$ws.Range("B15").Value = 0.9988676228301756
$ws.Range("C15").Value = 0.05670205453220944
$ws.Range("D15").Value = 0.5108987338151252
$ws.Range("E15").Value = 0.155996630156281
$ws.Range("G15").Value = 0.002533903603302967
$ws.Range("I15").Value = 1.4896555307906
$ws.Range("J15").Value = 0.05830828951633382
$ws.Range("K15").Value = 0.6843244697681428
$ws.Range("L15").Value = 0.445044134666162
$ws.Range("M15").Value = 0.3216889938347052
$ws.Range("N15").Value = 2.700046029964803
$ws.Range("O15").Value = 6.432856026956784
$ws.Range("B16").Value = 0.9735650647121759
$ws.Range("C16").Value = 0.05474332138027194
$ws.Range("D16").Value = 0.5089070364200978
$ws.Range("E16").Value = 0.1561635551253424
$ws.Range("G16").Value = 0.00253554127958792
$ws.Range("I16").Value = 1.493688891704856
$ws.Range("J16").Value = 0.05826320148381914
$ws.Range("K16").Value = 0.6597938920016873
$ws.Range("L16").Value = 0.4430662101676006
$ws.Range("M16").Value = 0.3168175206811767
$ws.Range("N16").Value = 2.714256769500707
$ws.Range("O16").Value = 6.443619095031295
$ws.Range("B17").Value = 0.9581301546212444
$ws.Range("C17").Value = 0.05353682094795431
$ws.Range("D17").Value = 0.5077364433427505
$ws.Range("E17").Value = 0.1562774494945494
$ws.Range("G17").Value = 0.002536568965898136
$ws.Range("I17").Value = 1.496306196607854
$ws.Range("J17").Value = 0.05823549705050457
$ws.Range("K17").Value = 0.6447959375355481
$ws.Range("L17").Value = 0.4418946931115215
$ws.Range("M17").Value = 0.313859997421396
$ws.Range("N17").Value = 2.723164698112146
$ws.Range("O17").Value = 6.450884982995092
$ws.Range("B18").Value = 0.9492847387001575
$ws.Range("C18").Value = 0.05284108353667705
$ws.Range("D18").Value = 0.5070820239174338
$ws.Range("E18").Value = 0.1563471896205311
$ws.Range("G18").Value = 0.002537168544573432
$ws.Range("I18").Value = 1.497864180215139
$ws.Range("J18").Value = 0.0582195460030599
$ws.Range("K18").Value = 0.6361883412735949
$ws.Range("L18").Value = 0.4412362945125636
$ws.Range("M18").Value = 0.3121703392057142
$ws.Range("N18").Value = 2.728358174068454
$ws.Range("O18").Value = 6.455307997838247
$ws.Range("B19").Value = 0.946295407545108
$ws.Range("C19").Value = 0.05260521165475041
$ws.Range("D19").Value = 0.5068636943888265
$ws.Range("E19").Value = 0.1563715298908637
$ws.Range("G19").Value = 0.00253737301001418
$ws.Range("I19").Value = 1.49840072086554
$ws.Range("J19").Value = 0.0582141424832221
$ws.Range("K19").Value = 0.6332772079335882
$ws.Range("L19").Value = 0.4410160252617175
$ws.Range("M19").Value = 0.3116002180215176
$ws.Range("N19").Value = 2.730128597179316
$ws.Range("O19").Value = 6.456847456702775
$ws.Range("B20").Value = 0.9597698854583427
$ws.Range("C20").Value = 0.05366544034211529
$ws.Range("D20").Value = 0.5078591024337697
$ws.Range("E20").Value = 0.1562648875243067
$ws.Range("G20").Value = 0.0025364586896691
$ws.Range("I20").Value = 1.496022139167763
$ws.Range("J20").Value = 0.05823844791320365
$ws.Range("K20").Value = 0.6463905504806746
$ws.Range("L20").Value = 0.4420178073177823
$ws.Range("M20").Value = 0.3141736488985103
$ws.Range("N20").Value = 2.722209202740157
$ws.Range("O20").Value = 6.450086281499608
$ws.Range("B21").Value = 1.005427083500649
$ws.Range("C21").Value = 0.05720624901914562
$ws.Range("D21").Value = 0.5114287615946864
$ws.Range("E21").Value = 0.1559570794886991
$ws.Range("G21").Value = 0.002533487890462763
$ws.Range("I21").Value = 1.488658365149945
$ws.Range("J21").Value = 0.0583199165390953
$ws.Range("K21").Value = 0.6906733322050798
$ws.Range("L21").Value = 0.4455677112710674
$ws.Range("M21").Value = 0.3229562453521027
$ws.Range("N21").Value = 2.696435860074409
$ws.Range("O21").Value = 6.430282067623693
$ws.Range("B22").Value = 1.03554996275426
$ws.Range("C22").Value = 0.05950435516469099
$ws.Range("D22").Value = 0.5139289467205828
$ws.Range("E22").Value = 0.1557934151360811
$ws.Range("G22").Value = 0.002531621417902773
$ws.Range("I22").Value = 1.484313106168955
$ws.Range("J22").Value = 0.0583730142483283
$ws.Range("K22").Value = 0.7197786292925059
$ws.Range("L22").Value = 0.4480243498043279
$ws.Range("M22").Value = 0.3287969003163767
$ws.Range("N22").Value = 2.680213482466563
$ws.Range("O22").Value = 6.419506599849484
$ws.Range("B23").Value = 1.019447123763229
$ws.Range("C23").Value = 0.05827928186630515
$ws.Range("D23").Value = 0.5125792936092495
$ws.Range("E23").Value = 0.1558773423604816
$ws.Range("G23").Value = 0.002532610735676186
$ws.Range("I23").Value = 1.486589525570118
$ws.Range("J23").Value = 0.05834468855342312
$ws.Range("K23").Value = 0.704229765335981
$ws.Range("L23").Value = 0.4467007382497314
$ws.Range("M23").Value = 0.3256704687371155
$ws.Range("N23").Value = 2.688814767550561
$ws.Range("O23").Value = 6.425059509862336
$ws.Range("B24").Value = 0.9590284748826718
$ws.Range("C24").Value = 0.05360729806854181
$ws.Range("D24").Value = 0.5078035904066098
$ws.Range("E24").Value = 0.156270553517535
$ws.Range("G24").Value = 0.002536508518318143
$ws.Range("I24").Value = 1.496150395624056
$ws.Range("J24").Value = 0.058237113901213
$ws.Range("K24").Value = 0.6456695794293523
$ws.Range("L24").Value = 0.4419621002729315
$ws.Range("M24").Value = 0.3140318138452471
$ws.Range("N24").Value = 2.722640957326405
$ws.Range("O24").Value = 6.450446608872568
$ws.Range("B25").Value = 0.8950959705092885
$ws.Range("C25").Value = 0.04849938725641323
$ws.Range("D25").Value = 0.5033743250440921
$ws.Range("E25").Value = 0.1568565598601808
$ws.Range("G25").Value = 0.002541037393988447
$ws.Range("I25").Value = 1.508474038000422
$ws.Range("J25").Value = 0.05812045935555599
$ws.Range("K25").Value = 0.5832248327562581
$ws.Range("L25").Value = 0.4374409668464949
$ws.Range("M25").Value = 0.3019152485003076
$ws.Range("N25").Value = 2.761799525235869
$ws.Range("O25").Value = 6.487156539939463
